$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.2556606666666667
$ws.Range("H2").Value2 = 0.7669820000000001
$ws.Range("I2").Value2 = 0.01354513404628681
$ws.Range("J2").Value2 = 0.01354513404628681
$ws.Range("M2").Value2 = 159.4836373333333
$ws.Range("N2").Value2 = 478.450912
$ws.Range("O2").Value2 = 0.2983285084902258
$ws.Range("P2").Value2 = 0.2983285084902258
$ws.Range("Q2").Value2 = 40.7736930430649
$ws.Range("R2").Value2 = 366.9632373875841
$ws.Range("S2").Value2 = 0.00404089963732892
$ws.Range("T2").Value2 = 0.00404089963732892
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.2556606666666667
$ws.Range("H3").Value2 = 0.7669820000000001
$ws.Range("I3").Value2 = 0.01354513404628681
$ws.Range("J3").Value2 = 0.01354513404628681
$ws.Range("O3").Value2 = 0.3227862111630279
$ws.Range("P3").Value2 = 0.3227862111630279
$ws.Range("Q3").Value2 = 44.11642038201801
$ws.Range("R3").Value2 = 397.047783438162
$ws.Range("S3").Value2 = 0.004372182498496252
$ws.Range("T3").Value2 = 0.004372182498496252
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.2556606666666667
$ws.Range("H4").Value2 = 0.7669820000000001
$ws.Range("I4").Value2 = 0.01354513404628681
$ws.Range("J4").Value2 = 0.01354513404628681
$ws.Range("M4").Value2 = 74.38770566666666
$ws.Range("N4").Value2 = 223.163117
$ws.Range("O4").Value2 = 0.1391489036280481
$ws.Range("P4").Value2 = 0.1391489036280482
$ws.Range("Q4").Value2 = 19.01801042254378
$ws.Range("R4").Value2 = 171.162093802894
$ws.Range("S4").Value2 = 0.001884790552035756
$ws.Range("T4").Value2 = 0.001884790552035757
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.2556606666666667
$ws.Range("H5").Value2 = 0.7669820000000001
$ws.Range("I5").Value2 = 0.01354513404628681
$ws.Range("J5").Value2 = 0.01354513404628681
$ws.Range("M5").Value2 = 58.41461433333333
$ws.Range("N5").Value2 = 175.243843
$ws.Range("O5").Value2 = 0.1092697975759847
$ws.Range("P5").Value2 = 0.1092697975759848
$ws.Range("Q5").Value2 = 14.93431924353622
$ws.Range("R5").Value2 = 134.408873191826
$ws.Range("S5").Value2 = 0.001480074055377339
$ws.Range("T5").Value2 = 0.001480074055377339
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.2556606666666667
$ws.Range("H6").Value2 = 0.7669820000000001
$ws.Range("I6").Value2 = 0.01354513404628681
$ws.Range("J6").Value2 = 0.01354513404628681
$ws.Range("M6").Value2 = 69.746216
$ws.Range("N6").Value2 = 209.238648
$ws.Range("O6").Value2 = 0.1304665791427133
$ws.Range("P6").Value2 = 0.1304665791427133
$ws.Range("Q6").Value2 = 17.83136408003734
$ws.Range("R6").Value2 = 160.482276720336
$ws.Range("S6").Value2 = 0.001767187303048538
$ws.Range("T6").Value2 = 0.001767187303048539
$ws.Range("I7").Value2 = 0.8835639662863414
$ws.Range("J7").Value2 = 0.8835639662863415
$ws.Range("M7").Value2 = 159.4836373333333
$ws.Range("N7").Value2 = 478.450912
$ws.Range("O7").Value2 = 0.2983285084902258
$ws.Range("P7").Value2 = 0.2983285084902258
$ws.Range("Q7").Value2 = 2659.71276638257
$ws.Range("R7").Value2 = 23937.41489744313
$ws.Range("S7").Value2 = 0.2635923202179123
$ws.Range("T7").Value2 = 0.2635923202179124
$ws.Range("I8").Value2 = 0.8835639662863414
$ws.Range("J8").Value2 = 0.8835639662863415
$ws.Range("O8").Value2 = 0.3227862111630279
$ws.Range("P8").Value2 = 0.3227862111630279
$ws.Range("S8").Value2 = 0.2852022649977455
$ws.Range("T8").Value2 = 0.2852022649977455
$ws.Range("I9").Value2 = 0.8835639662863414
$ws.Range("J9").Value2 = 0.8835639662863415
$ws.Range("M9").Value2 = 74.38770566666666
$ws.Range("N9").Value2 = 223.163117
$ws.Range("O9").Value2 = 0.1391489036280481
$ws.Range("P9").Value2 = 0.1391489036280482
$ws.Range("Q9").Value2 = 1240.565701483347
$ws.Range("R9").Value2 = 11165.09131335013
$ws.Range("S9").Value2 = 0.1229469571939941
$ws.Range("T9").Value2 = 0.1229469571939941
$ws.Range("I10").Value2 = 0.8835639662863414
$ws.Range("J10").Value2 = 0.8835639662863415
$ws.Range("M10").Value2 = 58.41461433333333
$ws.Range("N10").Value2 = 175.243843
$ws.Range("O10").Value2 = 0.1092697975759847
$ws.Range("P10").Value2 = 0.1092697975759848
$ws.Range("Q10").Value2 = 974.1820420169726
$ws.Range("R10").Value2 = 8767.638378152753
$ws.Range("S10").Value2 = 0.09654685574154273
$ws.Range("T10").Value2 = 0.09654685574154276
$ws.Range("I11").Value2 = 0.8835639662863414
$ws.Range("J11").Value2 = 0.8835639662863415
$ws.Range("M11").Value2 = 69.746216
$ws.Range("N11").Value2 = 209.238648
$ws.Range("O11").Value2 = 0.1304665791427133
$ws.Range("P11").Value2 = 0.1304665791427133
$ws.Range("Q11").Value2 = 1163.159457633616
$ws.Range("R11").Value2 = 10468.43511870254
$ws.Range("S11").Value2 = 0.1152755681351466
$ws.Range("T11").Value2 = 0.1152755681351467
$ws.Range("G12").Value2 = 1.942037333333333
$ws.Range("H12").Value2 = 5.826112
$ws.Range("I12").Value2 = 0.1028908996673717
$ws.Range("J12").Value2 = 0.1028908996673717
$ws.Range("M12").Value2 = 159.4836373333333
$ws.Range("N12").Value2 = 478.450912
$ws.Range("O12").Value2 = 0.2983285084902258
$ws.Range("P12").Value2 = 0.2983285084902258
$ws.Range("Q12").Value2 = 309.7231777571271
$ws.Range("R12").Value2 = 2787.508599814144
$ws.Range("S12").Value2 = 0.03069528863498448
$ws.Range("T12").Value2 = 0.03069528863498448
$ws.Range("G13").Value2 = 1.942037333333333
$ws.Range("H13").Value2 = 5.826112
$ws.Range("I13").Value2 = 0.1028908996673717
$ws.Range("J13").Value2 = 0.1028908996673717
$ws.Range("O13").Value2 = 0.3227862111630279
$ws.Range("P13").Value2 = 0.3227862111630279
$ws.Range("Q13").Value2 = 335.115043357888
$ws.Range("R13").Value2 = 3016.035390220992
$ws.Range("S13").Value2 = 0.03321176366678617
$ws.Range("T13").Value2 = 0.03321176366678617
$ws.Range("G14").Value2 = 1.942037333333333
$ws.Range("H14").Value2 = 5.826112
$ws.Range("I14").Value2 = 0.1028908996673717
$ws.Range("J14").Value2 = 0.1028908996673717
$ws.Range("M14").Value2 = 74.38770566666666
$ws.Range("N14").Value2 = 223.163117
$ws.Range("O14").Value2 = 0.1391489036280481
$ws.Range("P14").Value2 = 0.1391489036280482
$ws.Range("Q14").Value2 = 144.4637015456782
$ws.Range("R14").Value2 = 1300.173313911104
$ws.Range("S14").Value2 = 0.01431715588201828
$ws.Range("T14").Value2 = 0.01431715588201828
$ws.Range("G15").Value2 = 1.942037333333333
$ws.Range("H15").Value2 = 5.826112
$ws.Range("I15").Value2 = 0.1028908996673717
$ws.Range("J15").Value2 = 0.1028908996673717
$ws.Range("M15").Value2 = 58.41461433333333
$ws.Range("N15").Value2 = 175.243843
$ws.Range("O15").Value2 = 0.1092697975759847
$ws.Range("P15").Value2 = 0.1092697975759848
$ws.Range("Q15").Value2 = 113.4433618476018
$ws.Range("R15").Value2 = 1020.990256628416
$ws.Range("S15").Value2 = 0.01124286777906467
$ws.Range("T15").Value2 = 0.01124286777906467
$ws.Range("G16").Value2 = 1.942037333333333
$ws.Range("H16").Value2 = 5.826112
$ws.Range("I16").Value2 = 0.1028908996673717
$ws.Range("J16").Value2 = 0.1028908996673717
$ws.Range("M16").Value2 = 69.746216
$ws.Range("N16").Value2 = 209.238648
$ws.Range("O16").Value2 = 0.1304665791427133
$ws.Range("P16").Value2 = 0.1304665791427133
$ws.Range("Q16").Value2 = 135.4497553307307
$ws.Range("R16").Value2 = 1219.047797976576
$ws.Range("S16").Value2 = 0.01342382370451813
$ws.Range("T16").Value2 = 0.01342382370451813
